$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.034.34"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.270.56"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.87"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.14"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.530"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E7").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.44%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.80"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.112"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("E12").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.624.32"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.35"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.277.60"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.786"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.96%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.910.23"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.75"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.96%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0917"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.76%  "
$ws.Range("E20").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.22"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.02"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.60"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.94"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("E25").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.68"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("E28").ClearFormats()
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -9.79%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.02"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.27%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.89"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.35"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.96%  "
$ws.Range("E32").ClearFormats()
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0745"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.04"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.21"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.28%  "
$ws.Range("E36").ClearFormats()
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("E37").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.79"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("E40").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.23%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.74"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.016.58"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.56%  "
$ws.Range("E43").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +8.72%  "
$ws.Range("E44").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.38%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.36"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.54%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.91"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.11%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.29"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.87%  "
$ws.Range("E48").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.51"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.29%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.14"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.48%  "
$ws.Range("E51").ClearFormats()
